$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a duplicate of row 1 (same values/format), with H3 = 24
$ws.Range("A1:H1").Copy($ws.Range("A3:H3")) | Out-Null

# D3 needs the same hyperlink as D1 (mailto:javier@email.com)
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:javier@email.com") | Out-Null

# Re-apply row 1's formatting onto row 3 so D3 keeps the same style index as D1
# (Hyperlinks.Add() re-stamps the font on the target cell otherwise)
$ws.Range("A1:H1").Copy() | Out-Null
$ws.Range("A3:H3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New selection after the edit
$ws.Range("D10").Select() | Out-Null
